# The edit removes the original column A (a numeric "row index" style-1
# column whose values duplicated column F) and shifts the remaining
# columns B:F one place to the left to become A:E.
# Deleting the entire column A reproduces exactly this: B->A, C->B, D->C,
# E->D, F->E, with cell styles/values/shared-string references moving
# along with their cells, and the sheet dimension shrinking to A1:E7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A:A").Delete()
